$d = $word.ActiveDocument

# 1) "Assunto: ..." line -> strip the subject text and the "4" before "ª análise",
#    keep/extend the leading tabs.
$d.Content.Find.Execute(
    "Assunto: Retificação de área`t`t4ª análise",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "`t`t`t`tª análise", 2) | Out-Null

# 2) "Contribuinte: Carlos" -> "Contribuinte: " (clear the name, keep the label)
$d.Content.Find.Execute(
    "Contribuinte: Carlos",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Contribuinte: ", 2) | Out-Null

# 3) "Inscrição Imobiliária: 123456789123456" -> "Inscrição Imobiliária: "
$d.Content.Find.Execute(
    "Inscrição Imobiliária: 123456789123456",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Inscrição Imobiliária: ", 2) | Out-Null

# 4) "Endereço do imóvel: Rua Teste, nº 1 - bairro Teste, Itabira - MG"
#    -> "Endereço do imóvel: Rua , nº  - bairro ,  - " (clear the filled-in values)
$d.Content.Find.Execute(
    "Endereço do imóvel: Rua Teste, nº 1 - bairro Teste, Itabira - MG",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Endereço do imóvel: Rua , nº  - bairro ,  - ", 2) | Out-Null

# 5) "Dados recebidos: - Planta do imóvel / - Escritura /" -> "Dados recebidos: "
$d.Content.Find.Execute(
    "Dados recebidos: - Planta do imóvel`n- Escritura`n",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Dados recebidos: ", 2) | Out-Null

# 6) Insert a new blank paragraph right after "Dados recebidos: " and before the
#    closing "Após verificação..." paragraph.
$dadosPar = $d.Paragraphs(6)
$dadosPar.Range.InsertParagraphAfter() | Out-Null

# 7) Clear the final "Após verificação..." paragraph text entirely.
$lastPar = $d.Paragraphs($d.Paragraphs.Count)
$lastPar.Range.Find.Execute(
    "Após verificação dos arquivos apresentados à Prefeitura Municipal de Itabira referentes ao levantamento realizado, não foram identificados deslocamentos, sobreposições, nem invasão de vias públicas. Recomenda-se que a Prefeitura Municipal de Itabira opte pelo deferimento do processo XXXX/XX/XXXX.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "", 2) | Out-Null
